# Generate Report for Handback
#
# This script mirrors the "handback" report-generation commit:
#  - Overview sheet: status text for the two language columns changes
#    from "Ready for handoff" to "Handed back: in sync with en-US",
#    and those two columns are widened.
#  - zh-cn / de-de sheets: the "Latest Target File", "Latest Handback
#    File" and "Latest Handback DateTime" columns (J, K, L) are filled
#    in with the handback results (J becomes a hyperlink to the
#    original source file, like column A), and columns C, J, K are
#    widened to fit the new content.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/bf9bb1b3e7884d4065d951cdc4f33f1e0b60c474/e2e/"

# Target (handback) file names, in row order, shared by both language sheets.
$targetFiles = @(
    "93463dfa-68cf-422f-b64f-7fbecba66435.md",
    "963181bc-6958-4ff6-8c33-2b88c8d2e373.md",
    "abf76a13-d0ac-4ced-9749-2a9b9179f54b.md",
    "85857bf4-e956-4c4f-81ae-7a6552cbedc2.png",
    "0d76792b-9ff9-4a9f-a515-e06ca96e6eae.png"
)

# Latest Handback File names, per language.
$handbackFilesZhCn = @(
    "93463dfa-68cf-422f-b64f-7fbecba66435.e091aeb68b9cc4e11c6fd867a34e6865efe05053.zh-cn.xlf",
    "963181bc-6958-4ff6-8c33-2b88c8d2e373.d87e01e60e085c0193f3beaf6591cc8c1e2de774.zh-cn.xlf",
    "abf76a13-d0ac-4ced-9749-2a9b9179f54b.4fbf976bafc28d7afcf6b94380f60b52c6c3c1fa.zh-cn.xlf",
    "01aa2d6056b6c2f187da54934c0abd034583ae41.png",
    "7877a06446abf41082634a17e28121a3a47c84d5.png"
)
$handbackFilesDeDe = @(
    "93463dfa-68cf-422f-b64f-7fbecba66435.e091aeb68b9cc4e11c6fd867a34e6865efe05053.de-de.xlf",
    "963181bc-6958-4ff6-8c33-2b88c8d2e373.d87e01e60e085c0193f3beaf6591cc8c1e2de774.de-de.xlf",
    "abf76a13-d0ac-4ced-9749-2a9b9179f54b.4fbf976bafc28d7afcf6b94380f60b52c6c3c1fa.de-de.xlf",
    "01aa2d6056b6c2f187da54934c0abd034583ae41.png",
    "7877a06446abf41082634a17e28121a3a47c84d5.png"
)

$handbackDateZhCn = "2017-02-21 03:13:42"
$handbackDateDeDe = "2017-02-21 03:14:04"

# Column widths are stored internally in whole-pixel units, so the exact
# fractional "width" values from the original (non-Excel-authored) file
# can't be reproduced bit for bit; we request the nearest achievable
# ColumnWidth (raw width minus Excel's standard 0.8333... padding).
$pad = 5 / 6

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2:F6").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = 29.9777050018311 - $pad
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777050018311 - $pad

# ---------------------------------------------------------------------
# Helper to fill in one language sheet (zh-cn or de-de)
# ---------------------------------------------------------------------
function Update-LanguageSheet($ws, $handbackFiles, $handbackDate) {

    for ($i = 0; $i -lt 5; $i++) {
        $row = $i + 2

        # J: Latest Target File -> hyperlink to the same source file as column A
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 10), $baseUrl + $targetFiles[$i], "", "", $targetFiles[$i])

        # K: Latest Handback File
        $ws.Cells.Item($row, 11).Value = $handbackFiles[$i]

        # L: Latest Handback DateTime
        $ws.Cells.Item($row, 12).Value = $handbackDate
    }

    $ws.Columns.Item(3).ColumnWidth = 29.9777050018311 - $pad
    $ws.Columns.Item(10).ColumnWidth = 40 - $pad
    $ws.Columns.Item(11).ColumnWidth = 40 - $pad
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LanguageSheet $wsZhCn $handbackFilesZhCn $handbackDateZhCn

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LanguageSheet $wsDeDe $handbackFilesDeDe $handbackDateDeDe
